$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.530.66"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "1.879.22"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.09"
$ws.Range("E5").Value = "  -5.57%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4861"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2871"
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").Value = "1.877.28"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.62"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07229"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "88.03"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.967"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6587"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").Value = "30.487.39"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007775"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "2.120.98"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.721"
$ws.Range("E22").Value = "  -3.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "186.40"
$ws.Range("E23").Value = "  +6.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.002"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.206"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.83"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.833"
$ws.Range("E28").Value = "  -5.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.402"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.228"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08962"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.908"
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05187"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7327"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.703"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01809"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.649"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9153"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.016"
$ws.Range("E40").Value = "  -9.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4298"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.61"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9969"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.596"
$ws.Range("E44").Value = "  -6.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1328"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.196"
$ws.Range("E46").Value = "  -8.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05809"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.589"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.395"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3865"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.02"
$ws.Range("E51").Value = "  -1.26%  "
